# Fixed update to excel issue
# - Rename the "Requested quantity" header on the two existing sheets to
#   more descriptive, sheet-specific names.
# - Add a new "PO Forecast" sheet (after the existing sheets) with the
#   forecasted weekly PO quantities (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

# --- Rename the duplicated "Requested quantity" headers -------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page-margin defaults used by the rest of the workbook
# (0.75in / 1in / 0.5in == 54 / 72 / 36 points).
$wsForecast.PageSetup.LeftMargin   = 54
$wsForecast.PageSetup.RightMargin  = 54
$wsForecast.PageSetup.TopMargin    = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the header formatting (bold, border, centered) from the Weekly sheet
# so the new sheet matches the style used by the rest of the workbook.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows
$forecastRows = @(
    @(45452.99999999999, 206, 103.0797932281621, 309.2714091580665),
    @(45466.99999999999, 191, 78.35638421520532, 296.7790263754292),
    @(45473.99999999999, 183, 74.93544725018511, 290.2888150852232),
    @(45480.99999999999, 176, 68.70469146428412, 273.9389014587542),
    @(45487.99999999999, 168, 62.17738852583825, 272.0584060953273),
    @(45494.99999999999, 160, 56.04196733229229, 269.6610676512516),
    @(45522.99999999999, 130, 20.63457307624364, 236.484944716161),
    @(45529.99999999999, 122, 16.50443240007944, 233.3496141682121),
    @(45536.99999999999, 114, 11.01232093124245, 218.6853822594048),
    @(45543.99999999999, 107, 0.4960003169703763, 213.3758177835345),
    @(45550.99999999999, 99, -9.037897505229024, 209.7733976203321),
    @(45557.99999999999, 91, -13.66925073422756, 199.478488171575),
    @(45578.99999999999, 68, -32.81806061195083, 178.0026539615132),
    @(45585.99999999999, 61, -40.90415545273625, 170.9603959523931),
    @(45592.99999999999, 53, -56.12179264040296, 164.2743944942044),
    @(45599.99999999999, 46, -55.01872034198125, 147.7151835688697),
    @(45606.99999999999, 38, -72.66188504763245, 148.8637338525321),
    @(45634.99999999999, 7, -100.920708102669, 120.3476869394683),
    @(45641.99999999999, 0, -104.6924616242809, 107.1972519252471),
    @(45648.99999999999, 0, -117.2128987880874, 102.5377045910008),
    @(45655.99999999999, 0, -124.6684202716729, 92.4515448854842),
    @(45662.99999999999, 0, -126.3143224931872, 79.01227849346722),
    @(45669.99999999999, 0, -138.6505067023624, 75.24385728429395),
    @(45676.99999999999, 0, -143.5823729852044, 70.86835714937172),
    @(45683.99999999999, 0, -150.9904144627669, 66.9677586597146),
    @(45690.99999999999, 0, -156.2741661244176, 49.78541855380369),
    @(45697.99999999999, 0, -166.2402533440738, 47.79621960309726)
)

for ($i = 0; $i -lt $forecastRows.Length; $i++) {
    $row = $forecastRows[$i]
    $excelRow = $i + 2
    $wsForecast.Cells.Item($excelRow, 1).Value = $row[0]
    $wsForecast.Cells.Item($excelRow, 2).Value = $row[1]
    $wsForecast.Cells.Item($excelRow, 3).Value = $row[2]
    $wsForecast.Cells.Item($excelRow, 4).Value = $row[3]
}

# Copy the date-formatted column A style from the Weekly sheet onto the new
# sheet's "ds" column so it matches (custom date number format).
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A28").PasteSpecial(-4122)

Write-Host "PO Forecast sheet added with $($forecastRows.Length) data rows."
